# MassWateR ParameterMapping.xlsx update:
# "update parameter list with chl probe, phycocyanin, and phycoerythrin #67"
#
# - Row 31 (Chlorophyll / Chl a (probe)) and row 37 (Cyanobacteria / Cyanobacteria
#   (probe)) get an updated unit-of-measure string that now also lists "RFU".
# - Three new Cyanobacteria rows are inserted right after the existing
#   "Cyanobacteria (probe)" row (old row 37) and before the "Microcystins" row
#   (old row 39): Phycocyanin, Phycocyanin (probe), and Phycoerythrin.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUnits = "mg/l, ug/l, umol/l, ppm, RFU"

# Update the units-of-measure column for the two existing probe-related rows
# so that the allowed-units list also includes RFU.
$ws.Range("D31").Value = $newUnits
$ws.Range("D37").Value = $newUnits

# Insert three blank rows right before the old row 39 ("Microcystins"), so
# the Microcystins/Metals/Flow/Water Level/Air Temp rows shift down by three.
$ws.Rows.Item(39).Insert()
$ws.Rows.Item(39).Insert()
$ws.Rows.Item(39).Insert()

# Row 39: Phycocyanin
$ws.Range("A39").Value = "Cyanobacteria"
$ws.Range("B39").Value = "Phycocyanin"
$ws.Range("C39").Value = "Phycocyanin"
$ws.Range("D39").Value = "mg/l, ug/l, umol/l, ppm"
$ws.Range("B39").Style = $ws.Range("B37").Style

# Row 40: Phycocyanin (probe)  (note: WQX value keeps the upstream typo "Phcyocyanin (probe)")
$ws.Range("A40").Value = "Cyanobacteria"
$ws.Range("B40").Value = "Phycocyanin (probe)"
$ws.Range("C40").Value = "Phcyocyanin (probe)"
$ws.Range("D40").Value = $newUnits
$ws.Range("B40").Style = $ws.Range("B37").Style

# Row 41: Phycoerythrin
$ws.Range("A41").Value = "Cyanobacteria"
$ws.Range("B41").Value = "Phycoerythrin"
$ws.Range("C41").Value = "Phycoerythrin"
$ws.Range("D41").Value = $newUnits
$ws.Range("B41").Style = $ws.Range("B37").Style
